$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.861.54"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "'3.250.29"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'579.23"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'183.06"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "'3.821.10"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "'28.66"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").Value = "'67.893.12"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "'3.255.80"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'13.56"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "'379.44"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "'7.64"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'71.39"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").Value = "'0.513"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "'10.04"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").Value = "'22.88"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'7.02"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("D36").Value = "'162.41"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "'0.841"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'26.57"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.67"
$ws.Range("E41").Value = "  +4.33%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.60"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'25.54"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("D44").Value = "'347.29"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").Value = "'41.19"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'0.0688"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "'2.635.21"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").Value = "'0.0285"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "'0.103"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'0.992"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'6.19"
$ws.Range("E51").Value = "  +2.96%  "
